# First set of edits after R&R
#
# The "Rule" column of the hit_miss_rule table relabels the two
# "all-to-one-arm" baseline rows:
#   "All to control" -> "All to Status-quo"
#   "All to forcing"  -> "All to Structure"
#
# (Set C6 before C5 so the new shared-string entries land in the same
# order - "All to Structure" then "All to Status-quo" - as in the target
# workbook.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "All to Structure"
$ws.Range("C5").Value = "All to Status-quo"
